$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Stage the literal text "False" in a scratch cell (forced to text with a
# leading quote) so it can be copied as a *value* into D2/D3 without the
# copy itself picking up the "quote prefix" number format that a direct
# Range.Value assignment of "'False" would otherwise leave behind.
$scratch = $ws.Range("AZ1")
$scratch.Value = "'False"

# --- Row 2: replace with new purchase record (MRT1835) ---
$ws.Range("A2").Value = "MRT1835"
$ws.Range("B2").Value = "Purchase TL"
$scratch.Copy() | Out-Null
$ws.Range("D2").PasteSpecial(-4163) | Out-Null
$ws.Range("E2").Value = "Antonia Jane Allen and Mark Eamonn Reginald Flynn"
$ws.Range("I2").Value = " 23  Woodside     Leigh-on-Sea  Essex   SS9 4QX"
$ws.Range("K2").Value = "  9  Ennismore Gardens      Southend-on-Sea   Essex  SS2 5RA"
$ws.Range("L2").Value = "£355,000.00"
$ws.Range("M2").ClearContents()
$ws.Range("N2").ClearContents()
$ws.Range("R2").Value = "Paul Robinsons Solicitors"
$ws.Range("S2").Value = "470-474 London Road, Essex, SS0 9LD"
$ws.Range("T2").Value = "Andrew McClintock "
$ws.Range("V2").Value = "Individuals"
$ws.Range("AC2").Value = "Antonia Jane Allen"
$ws.Range("AD2").Value = "JE282389B"

# --- Row 3: replace with new purchase record (NBT1893) ---
$ws.Range("A3").Value = "NBT1893"
$ws.Range("B3").Value = "New Build Purchase"
$scratch.Copy() | Out-Null
$ws.Range("D3").PasteSpecial(-4163) | Out-Null
$ws.Range("E3").Value = "Fulvio Zaccagna and Maria Lyasheva"
$ws.Range("I3").Value = " 4  Grey Way     Cambridge  Cambridgeshire   CB5 8XT"
$ws.Range("K3").Value = "  Plot 9 Marleigh             "
$ws.Range("L3").Value = "£599,950.00"
$ws.Range("M3").ClearContents()
$ws.Range("N3").ClearContents()
$ws.Range("R3").Value = "Howard Kennedy"
$ws.Range("S3").Value = "1 London Bridge, SE1 9BG"
$ws.Range("T3").Value = " Plot 9 Marleigh  "
$ws.Range("V3").Value = "Individuals"
$ws.Range("AC3").Value = "Fulvio Zaccagna"
$ws.Range("AD3").Value = "ST077161B"

# Remove the scratch cell entirely (value + format) so it leaves no trace.
$scratch.Clear()

# --- Row 4: remove entirely (record no longer present) ---
$ws.Rows.Item(4).Delete()
